$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), matching the style of the existing
# header row (e.g. H1: bold font, thin border, centered alignment).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:I59 and J2:J59
$iVals = @(4,6,7,1,6,9,5,5,6,8,7,6,9,7,6,9,9,9,8,5,6,8,7,8,7,6,7,9,7,3,10,7,6,8,7,6,6,9,8,8,6,8,9,7,7,7,8,6,9,9,7,6,8,7,7,5,3,3)
$jVals = @(5,6,7,3,6,9,5,6,6,8,7,6,9,7,6,9,9,9,8,5,6,8,7,8,7,7,7,9,7,3,10,7,6,8,7,6,7,9,8,8,6,8,9,8,8,7,8,7,9,9,7,6,8,7,7,5,3,3)

for ($r = 0; $r -lt $iVals.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$r]
    $ws.Cells.Item($row, 10).Value = $jVals[$r]
}
